# Weekly update: insert a new price-report row at the top of the dated
# data block (row 346) for "Agrícola del Norte S.A. de Arica - Zanahoria".
# All rows that were previously 346..380 shift down to 347..381.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 346 - shifts existing rows 346:380 down to 347:381
$ws.Rows.Item(346).Insert()

# Populate the new row 346 with this week's data point
$ws.Range("A346").Value = 1
$ws.Range("B346").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C346").Value = "Arica y Parinacota"
$ws.Range("D346").Value = 44918
$ws.Range("E346").Value = 15
$ws.Range("F346").Value = 100114013
$ws.Range("G346").Value = "Zanahoria"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 210
$ws.Range("K346").Value = 22000
$ws.Range("L346").Value = 23000
$ws.Range("M346").Value = 22286
$ws.Range("N346").Value = '$/saco 25 kilos'
$ws.Range("O346").Value = "Región de Arica y Parinacota"
$ws.Range("P346").Value = 891
$ws.Range("Q346").Value = 25
$ws.Range("R346").Value = "Hortaliza"
